# "tidied up NJ tree plotting"
# The Colors lookup table's entry for "Outgroup" had its colour swatch
# changed from #ffff99 to #4d4d4d. Everything else in the workbook
# (the pop_names VLOOKUP results, the shared-string indices, etc.) just
# follows automatically from this single data edit when Excel recalculates.

$wb = $excel.ActiveWorkbook

$wsColors = $wb.Worksheets.Item("Colors")
$wsPop = $wb.Worksheets.Item("pop_names")

# Update the colour value used for the "Outgroup" row in the Colors table.
$wsColors.Range("B13").Value = "#4d4d4d"

# Force recalculation so the VLOOKUP formulas on pop_names pick up the
# new colour value.
$excel.Calculate()

# Restore the active cell / selection on each sheet.
$wsPop.Activate()
$wsPop.Range("E31").Select()

$wsColors.Activate()
$wsColors.Range("B16").Select()

# Leave pop_names as the active sheet/tab, matching the original workbook.
$wsPop.Activate()
